$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates per diff. Force text format ("@") before assigning so that
# numeric-looking strings (prices, percentages) are preserved as literal text,
# matching the original inlineStr storage instead of being auto-converted to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "278.89"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "6.82%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.31"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.11%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.820"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.64%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06266"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.78%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.860"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.63%"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.266"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.86%"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8792"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "3.13%"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9409"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.78%"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1448"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.34%"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.05137"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "6.24%"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07278"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.55%"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03164"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.83%"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09052"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.07%"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001548"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.88%"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "One"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006274"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.50%"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005860"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.61%"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.450"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.24%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.286"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "5.57%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1309"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.03%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.850"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-5.68%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.62%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.06%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004278"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001199"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.15%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "2.93%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04035"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.49%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006436"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "56.16%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1154"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "3.84%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002103"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-4.88%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01387"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.07%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005138"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.50%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.17%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.358"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "698.70%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-12.21%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002098"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.17%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001998"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.17%"

Write-Host "Applied 90 cell updates"
